$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# sportbet block
$ws.Range("A54").Value = "sportbet"
$ws.Range("A55").Value = "win  /div[2]/div[1]/div[1]/a[1]/span[2]"
$ws.Range("A56").Value = "draw  /div[2]/div[1]/div[2]/a[1]/span[2]"
$ws.Range("A57").Value = "lose  /div[2]/div[1]/div[3]/a[1]/span[2]"

# victorbet block
$ws.Range("A60").Value = "victorbet"
$ws.Range("A61").Value = "time   "
$ws.Range("B61").Value = " /td[1]"
$ws.Range("A62").Value = "teams"
$ws.Range("B62").Value = "/td[2]/a[1]"
$ws.Range("A63").Value = "league"
$ws.Range("B63").Value = "/td[2]/span[1]"
$ws.Range("A64").Value = "win"
$ws.Range("B64").Value = "/td[3]/span[1]/a[1]/span[1]"
$ws.Range("A65").Value = "draw"
$ws.Range("B65").Value = "/td[4]/span[1]/a[1]/span[1]"
$ws.Range("A66").Value = "lose"
$ws.Range("B66").Value = "/td[5]/span[1]/a[1]/span[1]"

# marathonbet block
$ws.Range("A68").Value = "marathonbet"
$ws.Range("A69").Value = "league"
$ws.Range("B69").Value = "div[1]/h2[1]"
$ws.Range("A71").Value = "start_time"
$ws.Range("B71").Value = "/tr[1]/td[1]/table[1]/tbody[1]/tr[1]/td[2]"
$ws.Range("B70").Value = "div[2]/div[1]/table[1]/tbody"
$ws.Range("A72").Value = "host"
$ws.Range("B72").Value = "/tr[1]/td[1]/table[1]/tbody[1]/tr[1]/td[1]/span[1]/div[1]"
$ws.Range("A73").Value = "client"
$ws.Range("B73").Value = "/tr[1]/td[1]/table[1]/tbody[1]/tr[1]/td[1]/span[1]/div[2]"
$ws.Range("A74").Value = "win"
$ws.Range("B74").Value = "/tr[1]/td[2]"
$ws.Range("A75").Value = "draw"
$ws.Range("B75").Value = "/tr[1]/td[3]"
$ws.Range("A76").Value = "lose"
$ws.Range("B76").Value = "/tr[1]/td[4]"

# coral block
$ws.Range("A80").Value = "coral"
$ws.Range("A81").Value = "win"
$ws.Range("B81").Value = "/div[5]/div[1]/span[2]"

# Column widths (engine quantizes ColumnWidth to whole display pixels at MDW=6;
# these are the closest achievable settings to the authored widths)
$ws.Columns(1).ColumnWidth = 13.666666666666666
$ws.Columns(2).ColumnWidth = 50.666666666666664

# Leave the view parked on the last cell entered, matching the authored selection
$ws.Range("B81").Select()
